$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the part number typo: LSH -> LTH
$ws.Range("D2").Value = "LTH-030-01-X-D-A-TR"

# Update the selected cell to match the saved cursor position
$ws.Range("G6").Select()
